$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (text columns)
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "nejnen"
$ws.Range("C2").Value = "sjnvjdn"
$ws.Range("D2").Value = "iejejiei"
# Force-text for numeric-looking / date-like strings so they keep their text type
$ws.Range("E2").Value = "'93837897"
$ws.Range("H2").Value = 88
$ws.Range("I2").Value = "'2020-02-02"
$ws.Range("J2").Value = "nwjnwjn"
$ws.Range("K2").Value = "il"
$ws.Range("M2").Value = "'89999"
$ws.Range("R2").Value = "'9999"

# Add new cells for row 2
$ws.Range("V2").Value = "2 Years"
$ws.Range("W2").Value = "None"

# Remove rows 3 and 4 entirely
$ws.Rows("3:4").Delete()
